$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 40
$ws.Range("D2").Value = 165.49
$ws.Range("H2").Value = 185.49
$ws.Range("I2").Value = 1.2366

$ws.Range("C4").Value = 10.70000000000002
$ws.Range("D4").Value = 154.1722
$ws.Range("H4").Value = 195.6388
$ws.Range("I4").Value = 0.978194
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 5
$ws.Range("P4").Value = 5

$ws.Range("C5").Value = 20
$ws.Range("D5").Value = 20
$ws.Range("H5").Value = 20
$ws.Range("I5").Value = 0.1333333333333333

$ws.Range("C6").Value = 0.3999999999999915
$ws.Range("D6").Value = 109.446
$ws.Range("H6").Value = 109.446
$ws.Range("I6").Value = 0.72964

$ws.Range("C7").Value = 1
$ws.Range("D7").Value = 2
$ws.Range("H7").Value = 2
$ws.Range("I7").Value = 0.01333333333333333

$ws.Range("C8").Value = 20.57000000000001
$ws.Range("D8").Value = 92.9925
$ws.Range("H8").Value = 92.9925
$ws.Range("I8").Value = 0.61995

$ws.Range("C9").Value = 9
$ws.Range("D9").Value = 165.6064
$ws.Range("H9").Value = 165.6064
$ws.Range("I9").Value = 1.104042666666667

$ws.Range("C10").Value = 8.099999999999994
$ws.Range("D10").Value = 92.3017
$ws.Range("H10").Value = 109.2017
$ws.Range("I10").Value = 0.7280113333333332

$ws.Range("C12").Value = 118.4160999999999
$ws.Range("D12").Value = 1162.3815
$ws.Range("E12").Value = 557.3342
$ws.Range("H12").Value = 1719.7157
$ws.Range("I12").Value = 1.146477133333333
$ws.Range("K12").Value = 4
$ws.Range("L12").Value = 83
$ws.Range("N12").Value = 38
$ws.Range("P12").Value = 123
$ws.Range("Q12").Value = 1.892307692307692
$ws.Range("S12").Value = 4
$ws.Range("T12").Value = 108
$ws.Range("W12").Value = 108
$ws.Range("X12").Value = 2.4

$ws.Range("C13").Value = 114
$ws.Range("E13").Value = 555.0622
$ws.Range("H13").Value = 733.8734
$ws.Range("I13").Value = 0.9784978666666666
$ws.Range("K13").Value = 10
$ws.Range("N13").Value = 36
$ws.Range("P13").Value = 145
$ws.Range("Q13").Value = 2.230769230769231
$ws.Range("S13").Value = 10
$ws.Range("T13").Value = 113
$ws.Range("W13").Value = 116
$ws.Range("X13").Value = 2.58

$ws.Range("C14").Value = 52.03400000000011
$ws.Range("D14").Value = 519.6352
$ws.Range("E14").Value = 268.7916
$ws.Range("H14").Value = 788.4268000000001
$ws.Range("I14").Value = 1.051235733333334

$ws.Range("C15").Value = 121.0000000000001
$ws.Range("D15").Value = 229.5
$ws.Range("E15").Value = 833.726
$ws.Range("H15").Value = 1063.226
$ws.Range("I15").Value = 1.635732307692308

$ws.Range("C16").Value = 54.00000000000011
$ws.Range("E16").Value = 579.1
$ws.Range("H16").Value = 747.9000000000001
$ws.Range("I16").Value = 1.4958
$ws.Range("K16").Value = 3
$ws.Range("N16").Value = 34
$ws.Range("P16").Value = 74
$ws.Range("Q16").Value = 1.233333333333333
$ws.Range("S16").Value = 3
$ws.Range("T16").Value = 48
$ws.Range("W16").Value = 50
$ws.Range("X16").Value = 1.25

$ws.Range("C17").Value = 69.60000000000002
$ws.Range("E17").Value = 718.2271
$ws.Range("H17").Value = 882.0271
$ws.Range("I17").Value = 1.7640542
$ws.Range("K17").Value = 2
$ws.Range("N17").Value = 26
$ws.Range("P17").Value = 61
$ws.Range("Q17").Value = 1.016666666666667
$ws.Range("S17").Value = 2
$ws.Range("T17").Value = 52
$ws.Range("W17").Value = 57
$ws.Range("X17").Value = 1.42

$ws.Range("C18").Value = 60
$ws.Range("E18").Value = 823.3088
$ws.Range("H18").Value = 851.4088
$ws.Range("I18").Value = 1.7028176
$ws.Range("K18").Value = 1
$ws.Range("N18").Value = 40
$ws.Range("P18").Value = 72
$ws.Range("Q18").Value = 1.2
$ws.Range("S18").Value = 1
$ws.Range("T18").Value = 66
$ws.Range("W18").Value = 66
$ws.Range("X18").Value = 1.65

